$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21 corresponds to month 2025-08-01 (A21 = 45870), last row of data.
$ws.Range("B21").Value = 6210
$ws.Range("C21").Value = 982
$ws.Range("D21").Value = 5590620
$ws.Range("E21").Value = 900.2608695652174
$ws.Range("F21").Value = 7.793785801076192
$ws.Range("G21").Value = 3.586497890295348
$ws.Range("H21").Value = 27.61140678370007
